$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the sampled "Media (velocidade - Rad/s)" readings in column D with
# a new batch of measurements. All dependent formulas (F6 average, J7, J11,
# J19) recompute automatically from these.
$ws.Range("D4").Value = 18.2605
$ws.Range("D5").Value = 17.73691
$ws.Range("D6").Value = 18.24741
$ws.Range("D7").Value = 18.05106
$ws.Range("D8").Value = 17.92016
$ws.Range("D9").Value = 17.98561
$ws.Range("D10").Value = 18.16887
$ws.Range("D11").Value = 18.41758
$ws.Range("D12").Value = 18.391407
$ws.Range("D13").Value = 17.95943
$ws.Range("D14").Value = 17.77617
$ws.Range("D15").Value = 17.94634
$ws.Range("D16").Value = 18.50921
$ws.Range("D17").Value = 17.448929
$ws.Range("D18").Value = 17.31802
$ws.Range("D19").Value = 17.31802
$ws.Range("D20").Value = 17.69763
$ws.Range("D21").Value = 16.91224
$ws.Range("D22").Value = 17.03005
$ws.Range("D23").Value = 16.32319
$ws.Range("D24").Value = 16.74207
$ws.Range("D25").Value = 16.37555
$ws.Range("D26").Value = 16.24465
$ws.Range("D27").Value = 16.3101
$ws.Range("D28").Value = 16.46718
$ws.Range("D29").Value = 17.095501
$ws.Range("D30").Value = 16.66353
$ws.Range("D31").Value = 16.85988
$ws.Range("D32").Value = 16.50645
$ws.Range("D33").Value = 16.12684

# Move the active selection to reflect where the editor ended up working.
$ws.Range("V25").Select()
